# Append: 2026-01-23 12:56 JST
# Update the "取得日時" (acquisition timestamp) column (A) for all data
# rows on the "ランサーズ" sheet from the old scrape time to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2026-01-23 12:42:09"
$newTimestamp = "2026-01-23 12:56:05"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
